# Updated cryptos list on Wed May  8 09:50:27 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "62.120.72"
$ws.Range("E2").Value = "  -3.27%  "

Set-TextValue "D3" "2.982.84"
$ws.Range("E3").Value = "  -4.37%  "

Set-TextValue "D5" "578.88"
$ws.Range("E5").Value = "  -2.74%  "

Set-TextValue "D6" "145.15"
$ws.Range("E6").Value = "  -8.24%  "

Set-TextValue "D8" "0.520"
$ws.Range("E8").Value = "  -4.13%  "

Set-TextValue "D9" "2.986.14"
$ws.Range("E9").Value = "  -4.31%  "

$ws.Range("E10").Value = "  -7.50%  "

$ws.Range("E11").Value = "  -4.82%  "

Set-TextValue "D12" "0.440"
$ws.Range("E12").Value = "  -2.94%  "

Set-TextValue "D13" "0.0000226"

Set-TextValue "D14" "34.43"
$ws.Range("E14").Value = "  -7.57%  "

$ws.Range("E15").Value = "  +1.45%  "

Set-TextValue "D16" "3.470.21"
$ws.Range("E16").Value = "  -4.48%  "

$ws.Range("E17").Value = "  -3.37%  "

Set-TextValue "D18" "62.131.16"
$ws.Range("E18").Value = "  -3.12%  "

Set-TextValue "D19" "2.984.79"
$ws.Range("E19").Value = "  -4.34%  "

Set-TextValue "D20" "454.79"
$ws.Range("E20").Value = "  -4.90%  "

$ws.Range("E21").Value = "  -4.98%  "

Set-TextValue "D22" "0.675"
$ws.Range("E22").Value = "  -5.81%  "

Set-TextValue "D23" "7.26"
$ws.Range("E23").Value = "  -4.28%  "

Set-TextValue "D24" "79.81"
$ws.Range("E24").Value = "  -1.99%  "

$ws.Range("E25").Value = "  -8.94%  "

Set-TextValue "D26" "12.14"
$ws.Range("E26").Value = "  -6.56%  "

Set-TextValue "D27" "0.999"
$ws.Range("E27").Value = "  -0.13%  "

$ws.Range("E28").Value = "  -6.26%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("E30").Value = "  -5.54%  "

$ws.Range("E31").Value = "  -4.06%  "

$ws.Range("E32").Value = "  -5.72%  "

$ws.Range("E33").Value = "  -2.52%  "

$ws.Range("E34").Value = "  -6.28%  "

Set-TextValue "D35" "1.02"
$ws.Range("E35").Value = "  -4.50%  "

Set-TextValue "D36" "0.0₃0777"
$ws.Range("E36").Value = "  -8.39%  "

$ws.Range("E37").Value = "  -5.75%  "

$ws.Range("E38").Value = "  -7.07%  "

Set-TextValue "D39" "49.91"
$ws.Range("E39").Value = "  -2.16%  "

Set-TextValue "D40" "8.95"
$ws.Range("E40").Value = "  -2.40%  "

$ws.Range("E41").Value = "  -12.97%  "

Set-TextValue "D42" "406.23"
$ws.Range("E42").Value = "  -9.60%  "

$ws.Range("E43").Value = "  -6.23%  "

$ws.Range("E44").Value = "  -1.77%  "

Set-TextValue "D45" "2.755.80"
$ws.Range("E45").Value = "  -2.73%  "

$ws.Range("E46").Value = "  -5.06%  "

Set-TextValue "D47" "38.20"
$ws.Range("E47").Value = "  -6.19%  "

Set-TextValue "D48" "127.13"
$ws.Range("E48").Value = "  -2.84%  "

$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("E50").Value = "  -2.95%  "

Set-TextValue "D51" "23.59"
$ws.Range("E51").Value = "  -9.15%  "
